$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 22.5
$ws.Range("J9").Value = 22.5
$ws.Range("L9").Value = 22.5
$ws.Range("N9").Value = -360.5

$ws.Range("H62").Value = 5009.625
$ws.Range("I62").Value = 3957.8462
$ws.Range("K62").Value = 3957.8462
$ws.Range("M62").Value = -3333.8462

$ws.Range("H65").Value = 5009.625
$ws.Range("I65").Value = 3957.8462
$ws.Range("K65").Value = 19789.231
$ws.Range("M65").Value = -16669.231

$ws.Range("H74").Value = 4624.3
$ws.Range("I74").Value = 3635.5454
$ws.Range("J74").Value = 5832.778
$ws.Range("K74").Value = 3635.5454
$ws.Range("L74").Value = 5832.778
$ws.Range("M74").Value = -2699.5454
$ws.Range("N74").Value = -7704.778

$ws.Range("H76").Value = 7521.7915
$ws.Range("I76").Value = 8271.308000000001
$ws.Range("K76").Value = 8271.308000000001
$ws.Range("M76").Value = -7956.308000000001

$ws.Range("H77").Value = 4624.3
$ws.Range("I77").Value = 3635.5454
$ws.Range("J77").Value = 5832.778
$ws.Range("K77").Value = 18177.727
$ws.Range("L77").Value = 29163.89
$ws.Range("M77").Value = -13497.727
$ws.Range("N77").Value = -38523.89

$ws.Range("H79").Value = 7521.7915
$ws.Range("I79").Value = 8271.308000000001
$ws.Range("K79").Value = 8271.308000000001
$ws.Range("M79").Value = -7179.308000000001

$ws.Range("H86").Value = 5068.6665
$ws.Range("I86").Value = 4456.2
$ws.Range("K86").Value = 4456.2
$ws.Range("M86").Value = -3333.2

$ws.Range("H87").Value = 69151.86
$ws.Range("J87").Value = 69151.86
$ws.Range("L87").Value = 69151.86
$ws.Range("N87").Value = -71647.86

$ws.Range("H89").Value = 5068.6665
$ws.Range("I89").Value = 4456.2
$ws.Range("K89").Value = 22281
$ws.Range("M89").Value = -16665

$ws.Range("H90").Value = 69151.86
$ws.Range("J90").Value = 69151.86
$ws.Range("L90").Value = 207455.58
$ws.Range("N90").Value = -219935.58

$ws.Range("H137").Value = 184234.95
$ws.Range("I137").Value = 244687.95
$ws.Range("J137").Value = 7194
$ws.Range("K137").Value = 734063.8500000001
$ws.Range("L137").Value = 21582
$ws.Range("M137").Value = -731513.8500000001
$ws.Range("N137").Value = -26682

$ws.Range("H138").Value = 2565.8933
$ws.Range("J138").Value = 2705.394
$ws.Range("L138").Value = 8116.181999999999
$ws.Range("N138").Value = -18396.182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 726.5
$ws.Range("I2").Value = 726.5
$ws.Range("K2").Value = 726.5
$ws.Range("M2").Value = -613.5

$ws.Range("H32").Value = 7048767
$ws.Range("I32").Value = 9093289
$ws.Range("K32").Value = 9093289
$ws.Range("M32").Value = -9093002

$ws.Range("H61").Value = 46881348
$ws.Range("I61").Value = 55559484
$ws.Range("J61").Value = 35723744
$ws.Range("K61").Value = 55559484
$ws.Range("L61").Value = 35723744
$ws.Range("M61").Value = -55559272
$ws.Range("N61").Value = -35724168

$ws.Range("H74").Value = 7820233.5
$ws.Range("I74").Value = 11365866
$ws.Range("J74").Value = 19841.3
$ws.Range("K74").Value = 11365866
$ws.Range("L74").Value = 19841.3
$ws.Range("M74").Value = -11364992
$ws.Range("N74").Value = -21589.3

$ws.Range("H77").Value = 7820233.5
$ws.Range("I77").Value = 11365866
$ws.Range("J77").Value = 19841.3
$ws.Range("K77").Value = 56829330
$ws.Range("L77").Value = 99206.5
$ws.Range("M77").Value = -56824962
$ws.Range("N77").Value = -107942.5

$ws.Range("H97").Value = 1167.9333
$ws.Range("I97").Value = 1034.3334
$ws.Range("J97").Value = 2370.3333
$ws.Range("K97").Value = 1034.3334
$ws.Range("L97").Value = 2370.3333
$ws.Range("M97").Value = -538.3334
$ws.Range("N97").Value = -3362.3333

$ws.Range("H116").Value = 726.5
$ws.Range("I116").Value = 726.5
$ws.Range("K116").Value = 726.5
$ws.Range("M116").Value = 1567.5

$ws.Range("H122").Value = 2901.8125
$ws.Range("I122").Value = 1955
$ws.Range("J122").Value = 4119.143
$ws.Range("K122").Value = 5865
$ws.Range("L122").Value = 12357.429
$ws.Range("M122").Value = -3415
$ws.Range("N122").Value = -17257.429

$ws.Range("H132").Value = 4985.3057
$ws.Range("I132").Value = 2539.0715
$ws.Range("J132").Value = 13547.125
$ws.Range("K132").Value = 7617.2145
$ws.Range("L132").Value = 40641.375
$ws.Range("M132").Value = -5087.2145
$ws.Range("N132").Value = -45701.375

$ws.Range("H136").Value = 46881348
$ws.Range("I136").Value = 55559484
$ws.Range("J136").Value = 35723744
$ws.Range("K136").Value = 166678452
$ws.Range("L136").Value = 107171232
$ws.Range("M136").Value = -166675902
$ws.Range("N136").Value = -107176332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 726.5
$ws.Range("I3").Value = 726.5
$ws.Range("K3").Value = 726.5
$ws.Range("M3").Value = -612.5

$ws.Range("H86").Value = 3655.5715
$ws.Range("I86").Value = 3458.1333
$ws.Range("J86").Value = 4149.1665
$ws.Range("K86").Value = 3458.1333
$ws.Range("L86").Value = 4149.1665
$ws.Range("M86").Value = -2335.1333
$ws.Range("N86").Value = -6395.1665

$ws.Range("H89").Value = 3655.5715
$ws.Range("I89").Value = 3458.1333
$ws.Range("J89").Value = 4149.1665
$ws.Range("K89").Value = 17290.6665
$ws.Range("L89").Value = 20745.8325
$ws.Range("M89").Value = -11674.6665
$ws.Range("N89").Value = -31977.8325

$ws.Range("H94").Value = 2121.5715
$ws.Range("J94").Value = 2888
$ws.Range("L94").Value = 2888
$ws.Range("N94").Value = -3790

$ws.Range("H107").Value = 565.5
$ws.Range("I107").Value = 565.5
$ws.Range("K107").Value = 565.5
$ws.Range("M107").Value = 1354.5

$ws.Range("H134").Value = 241840.78
$ws.Range("I134").Value = 1487.7941
$ws.Range("K134").Value = 4463.3823
$ws.Range("M134").Value = -1928.3823

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 594616.9399999999
$ws.Range("I31").Value = 10106.389
$ws.Range("K31").Value = 10106.389
$ws.Range("M31").Value = -9811.388999999999

$ws.Range("H34").Value = 594616.9399999999
$ws.Range("I34").Value = 10106.389
$ws.Range("K34").Value = 10106.389
$ws.Range("M34").Value = -9904.388999999999

$ws.Range("H58").Value = 10979
$ws.Range("I58").Value = 20013
$ws.Range("J58").Value = 7967.6665
$ws.Range("K58").Value = 20013
$ws.Range("L58").Value = 7967.6665
$ws.Range("M58").Value = -19810
$ws.Range("N58").Value = -8373.666499999999

$ws.Range("H122").Value = 1948.7
$ws.Range("I122").Value = 2081.9092
$ws.Range("K122").Value = 6245.7276
$ws.Range("M122").Value = -3795.7276

$ws.Range("H134").Value = 3360.1333
$ws.Range("I134").Value = 1400.5454
$ws.Range("K134").Value = 4201.6362
$ws.Range("M134").Value = -1666.6362

$ws.Range("H136").Value = 10979
$ws.Range("I136").Value = 20013
$ws.Range("J136").Value = 7967.6665
$ws.Range("K136").Value = 60039
$ws.Range("L136").Value = 23902.9995
$ws.Range("M136").Value = -57489
$ws.Range("N136").Value = -29002.9995

$ws.Range("H137").Value = 93959.336
$ws.Range("J137").Value = 93959.336
$ws.Range("L137").Value = 93959.336
$ws.Range("N137").Value = -104159.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 72782.71000000001
$ws.Range("J37").Value = 72782.71000000001
$ws.Range("L37").Value = 218348.13
$ws.Range("N37").Value = -218572.13

$ws.Range("H40").Value = 258.57144
$ws.Range("I40").Value = 52.75
$ws.Range("J40").Value = 533
$ws.Range("K40").Value = 211
$ws.Range("L40").Value = 2132
$ws.Range("M40").Value = -142
$ws.Range("N40").Value = -2270

$ws.Range("H127").Value = 1518.8
$ws.Range("J127").Value = 1518.8
$ws.Range("L127").Value = 4556.4
$ws.Range("N127").Value = -14476.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 41669656
$ws.Range("I132").Value = 41669656
$ws.Range("K132").Value = 125008968
$ws.Range("M132").Value = -125006438

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1275.6086
$ws.Range("I61").Value = 1003.375
$ws.Range("K61").Value = 1003.375
$ws.Range("M61").Value = -801.375

$ws.Range("H93").Value = 35716900
$ws.Range("I93").Value = 111113016
$ws.Range("J93").Value = 2948.8948
$ws.Range("K93").Value = 111113016
$ws.Range("L93").Value = 2948.8948
$ws.Range("M93").Value = -111111768
$ws.Range("N93").Value = -5444.8948

$ws.Range("H113").Value = 1275.6086
$ws.Range("I113").Value = 1003.375
$ws.Range("K113").Value = 1003.375
$ws.Range("M113").Value = 1166.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 439.7097
$ws.Range("I113").Value = 413.13043
$ws.Range("J113").Value = 516.125
$ws.Range("K113").Value = 1239.39129
$ws.Range("L113").Value = 1548.375
$ws.Range("M113").Value = 930.60871
$ws.Range("N113").Value = -5888.375

$ws.Range("H122").Value = 2863.6875
$ws.Range("I122").Value = 2401.138
$ws.Range("K122").Value = 7203.414
$ws.Range("M122").Value = -4753.414

$ws.Range("H132").Value = 253513.6
$ws.Range("I132").Value = 3244.6562
$ws.Range("J132").Value = 1254589.4
$ws.Range("K132").Value = 9733.9686
$ws.Range("L132").Value = 3763768.2
$ws.Range("M132").Value = -7203.9686
$ws.Range("N132").Value = -3768828.2
